$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells: Wins / Losses / Ties (columns AD, AE, AF) ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, centered, thin-border) from the
# existing header row (A1) onto the three new header cells.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Season record values for every player row (2-47) ---
$ws.Range("AD2:AD47").Value = 95
$ws.Range("AE2:AE47").Value = 67
$ws.Range("AF2:AF47").Value = 0

Write-Output "Season record columns (Wins/Losses/Ties) added."
